# Update "想去人数" (want-to-go count, column F) figures on each sheet
# to match freshly re-scraped totals (gh-pages data refresh at 456a3b4).
# Sheet order in workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型 (combined view)

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 273
$ws.Cells.Item(4, 6).Value = 1109
$ws.Cells.Item(5, 6).Value = 2674
$ws.Cells.Item(6, 6).Value = 228
$ws.Cells.Item(7, 6).Value = 685
$ws.Cells.Item(8, 6).Value = 60
$ws.Cells.Item(9, 6).Value = 255
$ws.Cells.Item(10, 6).Value = 187
$ws.Cells.Item(12, 6).Value = 98
$ws.Cells.Item(13, 6).Value = 126
$ws.Cells.Item(14, 6).Value = 1569
$ws.Cells.Item(17, 6).Value = 197
$ws.Cells.Item(18, 6).Value = 253

# Sheet 2: 演出 (Performances)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 26
$ws.Cells.Item(12, 6).Value = 46

# Sheet 3: 本地生活 (Local life)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 6351
$ws.Cells.Item(3, 6).Value = 796
$ws.Cells.Item(5, 6).Value = 252

# Sheet 4: 全部类型 (All types - combined view)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 6351
$ws.Cells.Item(3, 6).Value = 796
$ws.Cells.Item(5, 6).Value = 252
$ws.Cells.Item(9, 6).Value = 26
$ws.Cells.Item(11, 6).Value = 273
$ws.Cells.Item(12, 6).Value = 1109
$ws.Cells.Item(16, 6).Value = 2675
$ws.Cells.Item(18, 6).Value = 228
$ws.Cells.Item(21, 6).Value = 46
$ws.Cells.Item(22, 6).Value = 685
$ws.Cells.Item(23, 6).Value = 60
$ws.Cells.Item(24, 6).Value = 255
$ws.Cells.Item(26, 6).Value = 187
$ws.Cells.Item(28, 6).Value = 98
$ws.Cells.Item(29, 6).Value = 126
$ws.Cells.Item(31, 6).Value = 1569
$ws.Cells.Item(36, 6).Value = 197
$ws.Cells.Item(43, 6).Value = 253
